$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix Spanish axis name typos (extra "r" removed from "serorreversion")
$ws.Range("F10").Value = "Tasa de seroreversion"
$ws.Range("G10").Value = "Tasa de seroreversión Rhat"

# Extend the selection on the sheet to cover the full second table (A10:G14)
$ws.Range("A10:G14").Select()
